$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the bold/bordered
# header style already used by A1:H1 (copy style from H1, then overwrite
# the copied text with the real header labels).
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Column I (I0) is a constant 1 for every data row.
$ws.Range("I2:I18").Value = 1

# Column J (IF) mirrors column H (IP) for every data row.
$ws.Range("H2:H18").Copy($ws.Range("J2:J18"))
